$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 177303
$ws.Range("C4").Value = 167260
$ws.Range("C7").Value = 5.66
$ws.Range("C8").Value = 64.75
